# Updates the "cryptos" price/volume table (Price column D, Volume(1h) column E)
# with the latest scraped figures, matching the GitHub Actions refresh commit.
# Rows 27/28 also swap places (EthereumClassic <-> InjectiveProtocol).
#
# Note: several Price values look like plain numbers (e.g. "309.83", "18.00").
# Excel's COM layer auto-converts such strings to Doubles, which would silently
# drop significant trailing zeros / the original text formatting. Prefixing
# those with a leading apostrophe forces Excel to store them as text (exactly
# like the source inline-string cells), same as a user typing '309.83 by hand.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.028.12'
$ws.Range('E2').Value = '  +1.91%  '
$ws.Range('D3').Value = '2.287.70'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = "'309.83"
$ws.Range('E5').Value = '  +1.39%  '
$ws.Range('D6').Value = "'101.16"
$ws.Range('E6').Value = '  +3.79%  '
$ws.Range('D7').Value = "'0.530"
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = "'0.504"
$ws.Range('E9').Value = '  +2.48%  '
$ws.Range('D10').Value = "'36.27"
$ws.Range('E10').Value = '  +1.78%  '
$ws.Range('D11').Value = "'0.0819"
$ws.Range('E11').Value = '  +2.95%  '
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('D13').Value = "'6.86"
$ws.Range('E13').Value = '  +3.16%  '
$ws.Range('D14').Value = '2.633.90'
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('E15').Value = '  +2.65%  '
$ws.Range('D16').Value = '2.280.45'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = "'0.800"
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').Value = '42.913.88'
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('D19').Value = "'12.56"
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').Value = '0.0₃0915'
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('D21').Value = "'6.05"
$ws.Range('E21').Value = '  +1.07%  '
$ws.Range('D22').Value = "'67.76"
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = "'239.79"
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('D24').Value = "'2.62"
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('D25').Value = "'1.99"
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('D26').Value = "'0.998"
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'23.91"
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').Value = "'38.20"
$ws.Range('E28').Value = '  +2.59%  '
$ws.Range('D29').Value = "'9.65"
$ws.Range('E29').Value = '  +1.41%  '
$ws.Range('D30').Value = "'2.14"
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('D31').Value = "'166.20"
$ws.Range('E31').Value = '  +3.67%  '
$ws.Range('D32').Value = "'5.31"
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('D33').Value = "'1.00"
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').Value = "'18.00"
$ws.Range('E34').Value = '  +4.28%  '
$ws.Range('E35').Value = '  -2.24%  '
$ws.Range('D36').Value = "'0.0737"
$ws.Range('E36').Value = '  -0.40%  '
$ws.Range('E37').Value = '  +0.81%  '
$ws.Range('E38').Value = '  -0.75%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').Value = "'4.20"
$ws.Range('E41').Value = '  +3.20%  '
$ws.Range('D42').Value = "'2.29"
$ws.Range('E42').Value = '  -6.45%  '
$ws.Range('D43').Value = "'19.35"
$ws.Range('E43').Value = '  +2.46%  '
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('D45').Value = '1.951.85'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('D46').Value = "'2.99"
$ws.Range('E46').Value = '  +1.56%  '
$ws.Range('D47').Value = "'9.81"
$ws.Range('E47').Value = '  -1.42%  '
$ws.Range('D48').Value = "'54.82"
$ws.Range('E48').Value = '  +2.92%  '
$ws.Range('D49').Value = '2.503.13'
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('D50').Value = "'1.52"
$ws.Range('E50').Value = '  +0.24%  '
$ws.Range('D51').Value = "'72.64"
$ws.Range('E51').Value = '  +0.98%  '
